# "completed tank titrations 0216" - append the 2022-02-16 CRM titration
# accuracy check to the CRMAccuracyData log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$newRow = 59

$ws.Cells.Item($newRow, 1).Value = 20220216
$ws.Cells.Item($newRow, 2).Value = 2227.6280000000002
$ws.Cells.Item($newRow, 3).Value = 2224.4699999999998
$ws.Cells.Item($newRow, 4).Formula = "=100*(B" + $newRow + "-C" + $newRow + ")/C" + $newRow
$ws.Cells.Item($newRow, 5).Value = 180
$ws.Cells.Item($newRow, 6).Value = "CRM OPENED 20220118"

# Scroll the sheet down and leave the selection where the last entry lands,
# matching the author's on-screen state when they saved.
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Range("C59").Select()

$wb.Save()
